# Update crypto price/volume data per latest scrape
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "93.740.50"
Set-TextValue $ws.Range("E2") "  +3.68%  "

Set-TextValue $ws.Range("D3") "3.131.82"
Set-TextValue $ws.Range("E3") "  -1.04%  "

Set-TextValue $ws.Range("E4") "  +0.01%  "

Set-TextValue $ws.Range("D5") "244.04"
Set-TextValue $ws.Range("E5") "  +2.35%  "

Set-TextValue $ws.Range("D6") "618.13"
Set-TextValue $ws.Range("E6") "  -0.65%  "

Set-TextValue $ws.Range("D7") "1.10"
Set-TextValue $ws.Range("E7") "  -0.63%  "

Set-TextValue $ws.Range("D8") "0.415"
Set-TextValue $ws.Range("E8") "  +11.71%  "

Set-TextValue $ws.Range("E9") "  -0.10%  "

Set-TextValue $ws.Range("D10") "3.133.14"
Set-TextValue $ws.Range("E10") "  -0.93%  "

Set-TextValue $ws.Range("D11") "0.740"
Set-TextValue $ws.Range("E11") "  +0.03%  "

Set-TextValue $ws.Range("E12") "  -0.26%  "

Set-TextValue $ws.Range("D13") "0.0000259"
Set-TextValue $ws.Range("E13") "  +5.08%  "

Set-TextValue $ws.Range("D14") "34.91"
Set-TextValue $ws.Range("E14") "  -1.30%  "

Set-TextValue $ws.Range("D15") "93.051.98"
Set-TextValue $ws.Range("E15") "  +2.80%  "

Set-TextValue $ws.Range("D16") "5.51"
Set-TextValue $ws.Range("E16") "  -0.47%  "

Set-TextValue $ws.Range("D17") "3.711.81"

Set-TextValue $ws.Range("D18") "3.129.54"
Set-TextValue $ws.Range("E18") "  -1.68%  "

Set-TextValue $ws.Range("D19") "3.81"
Set-TextValue $ws.Range("E19") "  +3.18%  "

Set-TextValue $ws.Range("D20") "14.87"
Set-TextValue $ws.Range("E20") "  -1.45%  "

Set-TextValue $ws.Range("D21") "0.0000211"
Set-TextValue $ws.Range("E21") "  +3.85%  "

Set-TextValue $ws.Range("D22") "5.84"
Set-TextValue $ws.Range("E22") "  -0.54%  "

Set-TextValue $ws.Range("D23") "9.48"
Set-TextValue $ws.Range("E23") "  +4.19%  "

Set-TextValue $ws.Range("D24") "453.24"
Set-TextValue $ws.Range("E24") "  +2.83%  "

Set-TextValue $ws.Range("D25") "5.87"
Set-TextValue $ws.Range("E25") "  -1.91%  "

Set-TextValue $ws.Range("D26") "87.84"
Set-TextValue $ws.Range("E26") "  -1.44%  "

Set-TextValue $ws.Range("D27") "11.95"
Set-TextValue $ws.Range("E27") "  -0.29%  "

Set-TextValue $ws.Range("D28") "3.294.56"
Set-TextValue $ws.Range("E28") "  -1.06%  "

Set-TextValue $ws.Range("E29") "  -0.09%  "

Set-TextValue $ws.Range("E30") "  +6.65%  "

Set-TextValue $ws.Range("E31") "  -0.01%  "

Set-TextValue $ws.Range("D32") "0.227"
Set-TextValue $ws.Range("E32") "  -1.72%  "

Set-TextValue $ws.Range("D33") "9.29"
Set-TextValue $ws.Range("E33") "  -1.85%  "

Set-TextValue $ws.Range("E34") "  +0.20%  "

Set-TextValue $ws.Range("D35") "8.17"
Set-TextValue $ws.Range("E35") "  +4.07%  "

Set-TextValue $ws.Range("E36") "  -3.20%  "

Set-TextValue $ws.Range("D37") "26.39"
Set-TextValue $ws.Range("E37") "  +0.12%  "

Set-TextValue $ws.Range("B38") "PancakeSwap"
Set-TextValue $ws.Range("C38") "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
Set-TextValue $ws.Range("D38") "1.92"
Set-TextValue $ws.Range("E38") "  -1.68%  "

Set-TextValue $ws.Range("B39") "MantraDAO"
Set-TextValue $ws.Range("C39") "https://coinranking.com/coin/cTdD8lD-6+mantradao-om"
Set-TextValue $ws.Range("D39") "3.93"
Set-TextValue $ws.Range("E39") "  +5.07%  "

Set-TextValue $ws.Range("D40") "484.93"
Set-TextValue $ws.Range("E40") "  -4.00%  "

Set-TextValue $ws.Range("E41") "  -2.85%  "

Set-TextValue $ws.Range("D42") "3.52"
Set-TextValue $ws.Range("E42") "  +2.88%  "

Set-TextValue $ws.Range("D43") "0.439"
Set-TextValue $ws.Range("E43") "  -2.40%  "

Set-TextValue $ws.Range("D44") "23.09"
Set-TextValue $ws.Range("E44") "  +4.42%  "

Set-TextValue $ws.Range("E45") "  +0.03%  "

Set-TextValue $ws.Range("D46") "162.41"
Set-TextValue $ws.Range("E46") "  +3.29%  "

Set-TextValue $ws.Range("D47") "1.95"
Set-TextValue $ws.Range("E47") "  +1.39%  "

Set-TextValue $ws.Range("D48") "0.699"
Set-TextValue $ws.Range("E48") "  -3.30%  "

Set-TextValue $ws.Range("D49") "1.41"
Set-TextValue $ws.Range("E49") "  +2.70%  "

Set-TextValue $ws.Range("D50") "0.0334"
Set-TextValue $ws.Range("E50") "  +3.56%  "

Set-TextValue $ws.Range("D51") "4.48"
Set-TextValue $ws.Range("E51") "  +1.42%  "
